$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 33, shifting existing rows 33..141 down to 34..142.
$ws.Rows.Item(33).Insert()

# The cells that used to be in row 33 are now in row 34 (already shifted by Insert).
# Populate the new row 33 by copying the (now shifted) row 34 values, then override
# the date/volume/price columns with the new record's values.
$cols = @("A","B","C","E","F","G","H","I","N","O","Q","R")
foreach ($col in $cols) {
    $src = $ws.Range($col + "34")
    $dst = $ws.Range($col + "33")
    $dst.Value = $src.Value2
}

# Copy the date cell's number formatting/style from row 34 before setting new value.
$ws.Range("D34").Copy()
$ws.Range("D33").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("D33").Value = 44453
$ws.Range("J33").Value = 540
$ws.Range("K33").Value = 15000
$ws.Range("L33").Value = 16000
$ws.Range("M33").Value = 15500
$ws.Range("P33").Value = 1550
